# Scheduled market-data refresh: update Crafting Log profit-calc columns
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per latest Universalis pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 859.5454999999999
$ws.Range("I32").Value = 595
$ws.Range("J32").Value = 1080
$ws.Range("K32").Value = 595
$ws.Range("L32").Value = 1080
$ws.Range("M32").Value = -269
$ws.Range("N32").Value = -1732

$ws.Range("H62").Value = 4342.0527
$ws.Range("I62").Value = 4612.1875
$ws.Range("K62").Value = 4612.1875
$ws.Range("M62").Value = -3988.1875

$ws.Range("H65").Value = 4342.0527
$ws.Range("I65").Value = 4612.1875
$ws.Range("K65").Value = 23060.9375
$ws.Range("M65").Value = -19940.9375

$ws.Range("H106").Value = 6331.6665
$ws.Range("I106").Value = 6897.5
$ws.Range("K106").Value = 6897.5
$ws.Range("M106").Value = -6266.5

$ws.Range("H111").Value = 1484.4
$ws.Range("I111").Value = 1056.125
$ws.Range("J111").Value = 3197.5
$ws.Range("K111").Value = 3168.375
$ws.Range("L111").Value = 9592.5
$ws.Range("M111").Value = -101.375
$ws.Range("N111").Value = -15726.5

$ws.Range("H113").Value = 4383
$ws.Range("I113").Value = 3259.6
$ws.Range("K113").Value = 3259.6
$ws.Range("M113").Value = -5.599999999999909

$ws.Range("H132").Value = 41671692
$ws.Range("I132").Value = 47623836
$ws.Range("K132").Value = 142871508
$ws.Range("M132").Value = -142868978

$ws.Range("H137").Value = 2432.4614
$ws.Range("I137").Value = 1942.5
$ws.Range("J137").Value = 3216.4
$ws.Range("K137").Value = 5827.5
$ws.Range("L137").Value = 9649.200000000001
$ws.Range("M137").Value = -3277.5
$ws.Range("N137").Value = -14749.2

$ws.Range("H138").Value = 1162.7
$ws.Range("I138").Value = 864.1111
$ws.Range("K138").Value = 2592.3333
$ws.Range("M138").Value = 2547.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1055.7142
$ws.Range("I61").Value = 899.1667
$ws.Range("K61").Value = 899.1667
$ws.Range("M61").Value = -687.1667

$ws.Range("H132").Value = 8421.571
$ws.Range("I132").Value = 8787.4
$ws.Range("K132").Value = 26362.2
$ws.Range("M132").Value = -23832.2

$ws.Range("H136").Value = 1055.7142
$ws.Range("I136").Value = 899.1667
$ws.Range("K136").Value = 2697.5001
$ws.Range("M136").Value = -147.5001000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2564
$ws.Range("I107").Value = 1346
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 1346
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = 574
$ws.Range("N107").Value = -8840

$ws.Range("H134").Value = 6149.7334
$ws.Range("I134").Value = 6260.4287
$ws.Range("K134").Value = 18781.2861
$ws.Range("M134").Value = -16246.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2289.8
$ws.Range("I132").Value = 2289.8
$ws.Range("K132").Value = 6869.400000000001
$ws.Range("M132").Value = -4339.400000000001

$ws.Range("H134").Value = 1956.7407
$ws.Range("I134").Value = 1351.5416
$ws.Range("J134").Value = 6798.3335
$ws.Range("K134").Value = 4054.6248
$ws.Range("L134").Value = 20395.0005
$ws.Range("M134").Value = -1519.6248
$ws.Range("N134").Value = -25465.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H81").Value = 7780.8887
$ws.Range("J81").Value = 8128.5
$ws.Range("L81").Value = 24385.5
$ws.Range("N81").Value = -26631.5

$ws.Range("H84").Value = 7780.8887
$ws.Range("J84").Value = 8128.5
$ws.Range("L84").Value = 73156.5
$ws.Range("N84").Value = -84388.5

$ws.Range("H109").Value = 2548.625
$ws.Range("I109").Value = 898.1667
$ws.Range("K109").Value = 2694.5001
$ws.Range("M109").Value = -1654.5001

$ws.Range("H134").Value = 3125
$ws.Range("I134").Value = 5500
$ws.Range("J134").Value = 750
$ws.Range("K134").Value = 16500
$ws.Range("L134").Value = 2250
$ws.Range("M134").Value = -11430
$ws.Range("N134").Value = -12390

$ws.Range("H138").Value = 950.25
$ws.Range("I138").Value = 600.6667
$ws.Range("J138").Value = 1999
$ws.Range("K138").Value = 1802.0001
$ws.Range("L138").Value = 5997
$ws.Range("M138").Value = 3337.9999
$ws.Range("N138").Value = -16277

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1976
$ws.Range("I31").Value = 970
$ws.Range("K31").Value = 970
$ws.Range("M31").Value = -678

$ws.Range("H37").Value = 1976
$ws.Range("I37").Value = 970
$ws.Range("K37").Value = 970
$ws.Range("M37").Value = -693

$ws.Range("H80").Value = 6750
$ws.Range("I80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2002

$ws.Range("H83").Value = 6750
$ws.Range("I83").Value = 3000
$ws.Range("K83").Value = 15000
$ws.Range("M83").Value = -10008

$ws.Range("H113").Value = 1587.3334
$ws.Range("I113").Value = 1221.1428
$ws.Range("J113").Value = 2100
$ws.Range("K113").Value = 1221.1428
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 948.8571999999999
$ws.Range("N113").Value = -6440

$ws.Range("H132").Value = 900
$ws.Range("I132").Value = 900
$ws.Range("K132").Value = 2700
$ws.Range("M132").Value = -170

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 1000
$ws.Range("K32").Value = 1000
$ws.Range("M32").Value = -683

$ws.Range("H61").Value = 9571.857
$ws.Range("I61").Value = 8667.166999999999
$ws.Range("K61").Value = 8667.166999999999
$ws.Range("M61").Value = -8465.166999999999

$ws.Range("H113").Value = 9571.857
$ws.Range("I113").Value = 8667.166999999999
$ws.Range("K113").Value = 8667.166999999999
$ws.Range("M113").Value = -6497.166999999999

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 64000
$ws.Range("J16").Value = 64000
$ws.Range("L16").Value = 64000
$ws.Range("N16").Value = -64584

$ws.Range("H107").Value = 750
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -6840

$ws.Range("H123").Value = 69998.5
$ws.Range("J123").Value = 69998.5
$ws.Range("L123").Value = 69998.5
$ws.Range("N123").Value = -79798.5

$ws.Range("H132").Value = 1498.4
$ws.Range("I132").Value = 1498.4
$ws.Range("K132").Value = 4495.200000000001
$ws.Range("M132").Value = -1965.200000000001

$ws.Range("H136").Value = 2339.125
$ws.Range("I136").Value = 2184.8262
$ws.Range("J136").Value = 5888
$ws.Range("K136").Value = 6554.4786
$ws.Range("L136").Value = 17664
$ws.Range("M136").Value = -4004.4786
$ws.Range("N136").Value = -22764
